$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Label" in H1, copying the header format (bold, bordered, centered)
# from the neighboring header cell G1.
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new Label column: 0 for "Control" patient rows, 1 for "MDD" rows.
$controlRows = 2,3,4,5,6,12,13,14,15,16
$mddRows = 7,8,9,10,11,17,18,19,20,21

foreach ($r in $controlRows) {
    $ws.Cells.Item($r, 8).Value = 0
}
foreach ($r in $mddRows) {
    $ws.Cells.Item($r, 8).Value = 1
}
